$d = $word.ActiveDocument

$replacements = @(
    @("942÷3=314, 0", "433÷9=48, 1"),
    @("563÷4=140, 3", "674÷4=168, 2"),
    @("896÷5=179, 1", "366÷9=40, 6"),
    @("226÷3=75, 1", "844÷5=168, 4"),
    @("393÷4=98, 1", "866÷5=173, 1"),
    @("267÷5=53, 2", "580÷7=82, 6"),
    @("822÷4=205, 2", "746÷2=373, 0"),
    @("975÷7=139, 2", "402÷3=134, 0"),
    @("868÷4=217, 0", "405÷3=135, 0"),
    @("550÷9=61, 1", "768÷9=85, 3"),
    @("316÷7=45, 1", "119÷2=59, 1"),
    @("842÷3=280, 2", "476÷5=95, 1"),
    @("874÷5=174, 4", "752÷7=107, 3"),
    @("132÷5=26, 2", "905÷9=100, 5"),
    @("619÷8=77, 3", "904÷9=100, 4"),
    @("820÷8=102, 4", "384÷4=96, 0"),
    @("930÷3=310, 0", "259÷6=43, 1"),
    @("531÷4=132, 3", "636÷5=127, 1"),
    @("137÷4=34, 1", "889÷6=148, 1"),
    @("263÷2=131, 1", "945÷6=157, 3"),
    @("306÷2=153, 0", "519÷4=129, 3"),
    @("985÷4=246, 1", "321÷6=53, 3"),
    @("642÷2=321, 0", "810÷3=270, 0"),
    @("401÷9=44, 5", "529÷9=58, 7"),
    @("602÷3=200, 2", "306÷2=153, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
